# Appends the 2025-06-25 16:00 snapshot rows to each of the three portfolio
# sheets, mirroring the daily-refresh job that produced this commit
# ("自动添加每日生成文件 AIPEEarningYield.xlsx (Colab)"). Each sheet gets one
# new row per current holding (same columns as every existing block):
#   A 组合名称 | B 股票代码 | C 股票名称 | D 配置比例(%) | E 标的数量
#   F 价格 | G 价值 | H 总价值 | I 修改时间
# Column B holds codes like "000333"/"HK02899" that must stay text (leading
# zeros / letters), so NumberFormat is forced to "@" before the values land.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param($ws, [int]$startRow, $rows)

    $endRow = $startRow + $rows.Count - 1

    # A:C and I are text columns; lock the format before writing so codes
    # such as "000333" or "100000" aren't auto-coerced into numbers. D:H stay
    # on the sheet's default (General) numeric format, untouched.
    $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 3)).NumberFormat = "@"
    $ws.Range($ws.Cells.Item($startRow, 9), $ws.Cells.Item($endRow, 9)).NumberFormat = "@"

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $row = $rows[$i]
        for ($j = 0; $j -lt $row.Count; $j++) {
            $ws.Cells.Item($r, $j + 1).Value = $row[$j]
        }
    }
}

# --- 大智投资组合收益: rows 138-148 (11 holdings @ 202506251600) ---
$ws1 = $wb.Worksheets.Item("大智投资组合收益")
$rows1 = @(
    @("大智 (稳健智远)","000333","美的集团",3.06,42.89719527444501,73.18000000000001,3139.216750183886,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","510050","上证50ETF",5.1,1852.638869852698,2.83,5242.968001683135,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","510300","沪深300ETF",4.97,1276.881426590205,4,5107.52570636082,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","518880","黄金ETF",4.85,673.9096418114974,7.39,4980.192252986965,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","600085","同仁堂",2,56.06334956913254,36.63,2053.600494717325,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","600900","长江电力",20.11,669.0422707938175,30.88,20660.02532211308,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","600989","宝丰能源",4.86,308.5861230826723,16.17,4989.837610246812,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","601899","XD紫金矿",10.01,541.4194359293219,18.99,10281.55508829782,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","HK02899","紫金矿业",9.880000000000001,521.4694874201556,19.46,10147.79622519623,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","HK06881","中国银河",5.43,600.7968033130801,9.289999999999999,5581.402302778513,102752.6611071341,"202506251600"),
    @("大智 (稳健智远)","100000","现金",29.75,30568.54135256952,1,30568.54135256952,102752.6611071341,"202506251600")
)
Append-Rows $ws1 138 $rows1

# --- 大成投资组合收益: rows 93-99 (7 holdings @ 202506251600) ---
$ws2 = $wb.Worksheets.Item("大成投资组合收益")
$rows2 = @(
    @("大成 (锐进先锋)","000725","京东方A",4.94,1264.047242468624,3.96,5005.627080175751,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","159781","科创创业ETF",5.04,9277.629760760277,0.55,5102.696368418153,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","513100","纳指ETF",4.92,3131.938709046463,1.59,4979.782547383877,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","513290","纳指生物科技ETF",0.95,870.2909333102562,1.11,966.0229359743844,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","603119","浙江荣泰",47.53,1069.978093782073,44.98,48127.61465831765,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","688290","景业智能",7.48,147.121987895035,51.49,7575.311156715353,101259.9173862028,"202506251600"),
    @("大成 (锐进先锋)","100000","现金",29.14,29502.86263921768,1,29502.86263921768,101259.9173862028,"202506251600")
)
Append-Rows $ws2 93 $rows2

# --- 我的投资组合收益: rows 196-210 (15 holdings @ 202506251600) ---
$ws3 = $wb.Worksheets.Item("我的投资组合收益")
$rows3 = @(
    @("范式进化投资组合","000333","美的集团",1.02,14.02515882310653,73.18000000000001,1026.361122674936,100570.9575779558,"202506251600"),
    @("范式进化投资组合","000725","京东方A",5.06,1284.618017653691,3.96,5087.087349908616,100570.9575779558,"202506251600"),
    @("范式进化投资组合","159781","科创创业ETF",5.16,9428.611488061997,0.55,5185.736318434098,100570.9575779558,"202506251600"),
    @("范式进化投资组合","510050","上证50ETF",5.11,1817.150577699221,2.83,5142.536134888795,100570.9575779558,"202506251600"),
    @("范式进化投资组合","510300","沪深300ETF",4.98,1252.422077361618,4,5009.688309446472,100570.9575779558,"202506251600"),
    @("范式进化投资组合","513100","纳指ETF",1.01,636.5814125697908,1.59,1012.164445985967,100570.9575779558,"202506251600"),
    @("范式进化投资组合","513290","纳指生物科技ETF",0.98,884.4538210040457,1.11,981.7437413144909,100570.9575779558,"202506251600"),
    @("范式进化投资组合","518880","黄金ETF",0.97,132.2001081659486,7.39,976.95879934636,100570.9575779558,"202506251600"),
    @("范式进化投资组合","600085","同仁堂",1,27.49471300507762,36.63,1007.131337375993,100570.9575779558,"202506251600"),
    @("范式进化投资组合","600900","长江电力",1.01,32.81132034584936,30.88,1013.213572279828,100570.9575779558,"202506251600"),
    @("范式进化投资组合","600989","宝丰能源",4.87,302.6749902285196,16.17,4894.254591995163,100570.9575779558,"202506251600"),
    @("范式进化投资组合","601899","XD紫金矿",10.03,531.0482559694855,18.99,10084.60638086053,100570.9575779558,"202506251600"),
    @("范式进化投资组合","603119","浙江荣泰",1.08,24.16423640557475,44.98,1086.907353522752,100570.9575779558,"202506251600"),
    @("范式进化投资组合","HK06881","中国银河",1.09,117.8576436007749,9.289999999999999,1094.897509051199,100570.9575779558,"202506251600"),
    @("范式进化投资组合","100000","现金",56.64,56967.67061087058,1,56967.67061087058,100570.9575779558,"202506251600")
)
Append-Rows $ws3 196 $rows3
